$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Change 1: "Critério" paragraph run — split text before the formula,
# insert three manual line breaks, keep formula, existing trailing <w:br/> stays.
$find.Execute(
    "5,0 pontos.(Nota final+P_recuperação)/2",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "5,0 pontos.^l^l^l(Nota final+P_recuperação)/2",
    2
) | Out-Null

# Change 2: "Norma de recuperação" paragraph run — split text before the
# formula with two manual line breaks (no trailing break after).
$find.Execute(
    "e sua nota final.(Nota final+P_recuperação)/2",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "e sua nota final.^l^l(Nota final+P_recuperação)/2",
    2
) | Out-Null

# Change 3: Bibliography run — split between the two references with two
# manual line breaks.
$find.Execute(
    "1473p.B)OGA",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "1473p.^l^lB)OGA",
    2
) | Out-Null
